$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H19").Value = 1319.8334
$ws.Range("I19").Value = 1461.875
$ws.Range("J19").Value = 1206.2
$ws.Range("K19").Value = 1461.875
$ws.Range("L19").Value = 1206.2
$ws.Range("M19").Value = -1286.875
$ws.Range("N19").Value = -1556.2
$ws.Range("H28").Value = 1032.7858
$ws.Range("J28").Value = 1010.1429
$ws.Range("L28").Value = 1010.1429
$ws.Range("N28").Value = -1980.1429
$ws.Range("H55").Value = 97.25
$ws.Range("I55").Value = 113
$ws.Range("K55").Value = 113
$ws.Range("M55").Value = 101
$ws.Range("H107").Value = 1522.1
$ws.Range("I107").Value = 1176.0667
$ws.Range("J107").Value = 2560.2
$ws.Range("K107").Value = 1176.0667
$ws.Range("L107").Value = 2560.2
$ws.Range("M107").Value = 743.9332999999999
$ws.Range("N107").Value = -6400.2
$ws.Range("H113").Value = 5834.1665
$ws.Range("J113").Value = 5219.8
$ws.Range("L113").Value = 5219.8
$ws.Range("N113").Value = -11727.8
$ws.Range("H116").Value = 498621.47
$ws.Range("I116").Value = 617803.25
$ws.Range("J116").Value = 6996.625
$ws.Range("K116").Value = 617803.25
$ws.Range("L116").Value = 6996.625
$ws.Range("M116").Value = -614361.25
$ws.Range("N116").Value = -13880.625
$ws.Range("H127").Value = 1471.3
$ws.Range("I127").Value = 894.7857
$ws.Range("K127").Value = 2684.3571
$ws.Range("M127").Value = 2275.6429
$ws.Range("H137").Value = 40591.312
$ws.Range("I137").Value = 32168.6
$ws.Range("J137").Value = 54629.168
$ws.Range("K137").Value = 96505.79999999999
$ws.Range("L137").Value = 163887.504
$ws.Range("M137").Value = -93955.79999999999
$ws.Range("N137").Value = -168987.504

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6786
$ws.Range("I61").Value = 972.7368
$ws.Range("K61").Value = 972.7368
$ws.Range("M61").Value = -760.7368
$ws.Range("H63").Value = 3711.5293
$ws.Range("I63").Value = 1709.4
$ws.Range("K63").Value = 1709.4
$ws.Range("M63").Value = -1023.4
$ws.Range("H66").Value = 3711.5293
$ws.Range("I66").Value = 1709.4
$ws.Range("K66").Value = 8547
$ws.Range("M66").Value = -5115
$ws.Range("H110").Value = 38046.832
$ws.Range("I110").Value = 40798.59
$ws.Range("K110").Value = 40798.59
$ws.Range("M110").Value = -38753.59
$ws.Range("H122").Value = 3018.9
$ws.Range("I122").Value = 2909.889
$ws.Range("K122").Value = 8729.667000000001
$ws.Range("M122").Value = -6279.667000000001
$ws.Range("H132").Value = 3314.65
$ws.Range("I132").Value = 2853
$ws.Range("K132").Value = 8559
$ws.Range("M132").Value = -6029
$ws.Range("H136").Value = 6786
$ws.Range("I136").Value = 972.7368
$ws.Range("K136").Value = 2918.2104
$ws.Range("M136").Value = -368.2103999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1687
$ws.Range("I105").Value = 1492.56
$ws.Range("K105").Value = 1492.56
$ws.Range("M105").Value = 254.4400000000001
$ws.Range("H139").Value = 109999.8
$ws.Range("J139").Value = 109999.8
$ws.Range("L139").Value = 109999.8
$ws.Range("N139").Value = -120279.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5887692
$ws.Range("I31").Value = 12504743
$ws.Range("J31").Value = 5869.222
$ws.Range("K31").Value = 12504743
$ws.Range("L31").Value = 5869.222
$ws.Range("M31").Value = -12504448
$ws.Range("N31").Value = -6459.222
$ws.Range("H34").Value = 5887692
$ws.Range("I34").Value = 12504743
$ws.Range("J34").Value = 5869.222
$ws.Range("K34").Value = 12504743
$ws.Range("L34").Value = 5869.222
$ws.Range("M34").Value = -12504541
$ws.Range("N34").Value = -6273.222
$ws.Range("H58").Value = 1540.5
$ws.Range("I58").Value = 1376.04
$ws.Range("K58").Value = 1376.04
$ws.Range("M58").Value = -1173.04
$ws.Range("H96").Value = 5999.5
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 5999.5
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 5999.5
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -11491.5
$ws.Range("H99").Value = 2636.9285
$ws.Range("J99").Value = 4031.25
$ws.Range("L99").Value = 4031.25
$ws.Range("N99").Value = -7027.25
$ws.Range("H105").Value = 2625.4443
$ws.Range("I105").Value = 1526
$ws.Range("J105").Value = 3999.75
$ws.Range("K105").Value = 1526
$ws.Range("L105").Value = 3999.75
$ws.Range("M105").Value = 221
$ws.Range("N105").Value = -7493.75
$ws.Range("H116").Value = 39500
$ws.Range("J116").Value = 39500
$ws.Range("L116").Value = 39500
$ws.Range("N116").Value = -48678
$ws.Range("H126").Value = 2636.9285
$ws.Range("J126").Value = 4031.25
$ws.Range("L126").Value = 12093.75
$ws.Range("N126").Value = -17033.75
$ws.Range("H132").Value = 67792.60000000001
$ws.Range("I132").Value = 77914.62
$ws.Range("K132").Value = 233743.86
$ws.Range("M132").Value = -231213.86
$ws.Range("H134").Value = 2684.5
$ws.Range("I134").Value = 2121.6
$ws.Range("K134").Value = 6364.799999999999
$ws.Range("M134").Value = -3829.799999999999
$ws.Range("H136").Value = 1540.5
$ws.Range("I136").Value = 1376.04
$ws.Range("K136").Value = 4128.12
$ws.Range("M136").Value = -1578.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 550
$ws.Range("J113").Value = 400
$ws.Range("L113").Value = 1200
$ws.Range("N113").Value = -5540
$ws.Range("H114").Value = 714.5
$ws.Range("I114").Value = 714.5
$ws.Range("K114").Value = 2143.5
$ws.Range("M114").Value = 1110.5
$ws.Range("H122").Value = 662
$ws.Range("J122").Value = 666.2
$ws.Range("L122").Value = 5995.8
$ws.Range("N122").Value = -10895.8
$ws.Range("H124").Value = 7818.4585
$ws.Range("J124").Value = 9000
$ws.Range("L124").Value = 27000
$ws.Range("N124").Value = -36820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 42600
$ws.Range("J52").Value = 42600
$ws.Range("L52").Value = 42600
$ws.Range("N52").Value = -43118
$ws.Range("H102").Value = 3540.125
$ws.Range("I102").Value = 2664.6
$ws.Range("K102").Value = 2664.6
$ws.Range("M102").Value = -1042.6
$ws.Range("H113").Value = 1728.25
$ws.Range("I113").Value = 1666.6666
$ws.Range("K113").Value = 1666.6666
$ws.Range("M113").Value = 503.3334
$ws.Range("H118").Value = 23724.75
$ws.Range("J118").Value = 23724.75
$ws.Range("L118").Value = 23724.75
$ws.Range("N118").Value = -27038.75
$ws.Range("H122").Value = 3687.0715
$ws.Range("I122").Value = 3551.4614
$ws.Range("J122").Value = 5450
$ws.Range("K122").Value = 10654.3842
$ws.Range("L122").Value = 16350
$ws.Range("M122").Value = -8204.3842
$ws.Range("N122").Value = -21250

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3388
$ws.Range("I46").Value = 1020.4
$ws.Range("J46").Value = 5361
$ws.Range("K46").Value = 1020.4
$ws.Range("L46").Value = 5361
$ws.Range("M46").Value = -832.4
$ws.Range("N46").Value = -5737
$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -40990
$ws.Range("H99").Value = 56129
$ws.Range("I99").Value = 56129
$ws.Range("K99").Value = 56129
$ws.Range("M99").Value = -53134
$ws.Range("H100").Value = 2539
$ws.Range("I100").Value = 2420.7
$ws.Range("K100").Value = 2420.7
$ws.Range("M100").Value = -1879.7
$ws.Range("H122").Value = 3920.6
$ws.Range("I122").Value = 3920.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11761.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9311.799999999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 5299.6665
$ws.Range("I132").Value = 5449.5
$ws.Range("K132").Value = 16348.5
$ws.Range("M132").Value = -13818.5
$ws.Range("H136").Value = 6068.353
$ws.Range("I136").Value = 6995
$ws.Range("K136").Value = 20985
$ws.Range("M136").Value = -18435
$ws.Range("H140").Value = 107998
$ws.Range("J140").Value = 107998
$ws.Range("L140").Value = 107998
$ws.Range("N140").Value = -118358

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 749.7143
$ws.Range("I100").Value = 791.5
$ws.Range("J100").Value = 499
$ws.Range("K100").Value = 1583
$ws.Range("L100").Value = 998
$ws.Range("M100").Value = -1042
$ws.Range("N100").Value = -2080
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H117").Value = 75000
$ws.Range("J117").Value = 75000
$ws.Range("L117").Value = 75000
$ws.Range("N117").Value = -84178
$ws.Range("H132").Value = 26693.107
$ws.Range("I132").Value = 29456.28
$ws.Range("K132").Value = 88368.84
$ws.Range("M132").Value = -85838.84
$ws.Range("H136").Value = 28158.885
$ws.Range("J136").Value = 2638
$ws.Range("L136").Value = 7914
$ws.Range("N136").Value = -13014
$ws.Range("H139").Value = 109143
$ws.Range("J139").Value = 109143
$ws.Range("L139").Value = 109143
$ws.Range("N139").Value = -119423
